# Revert "Powerpoint writer: consolidate text run nodes."
#
# The caption text box on slide 1 ("The picture first") currently has
# its text consolidated into 3 runs: "The ", "picture ", "first".
# Split it back into 5 runs: "The", " ", "picture", " ", "first"
# (i.e. trailing spaces get pulled out into their own run), leaving
# the visible text and all formatting untouched.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Re-assigning each character range's Text to itself (same content)
# forces the writer to emit a separate run for that span, without
# altering the paragraph's overall text or formatting.
$tr.Characters(1, 3).Text  = "The"
$tr.Characters(4, 1).Text  = " "
$tr.Characters(5, 7).Text  = "picture"
$tr.Characters(12, 1).Text = " "
$tr.Characters(13, 5).Text = "first"
